$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 50.021050241203319
$ws.Range("C3").Value = 50.081785624074378

$ws.Range("B4").Value = 49.950097687778069
$ws.Range("C4").Value = 50.037720871183069

$ws.Range("B5").Value = 49.992580497477412
$ws.Range("C5").Value = 50.098066932655975

$ws.Range("B6").Value = 49.742032858979286
$ws.Range("C6").Value = 49.861034278436001

$ws.Range("B7").Value = 50.301427828224497
$ws.Range("C7").Value = 50.430793361179198

$ws.Range("B9").Value = 0.062179543649808135
$ws.Range("C9").Value = 0.057790335571873179

$ws.Range("B13").Value = 1.4502498382035611
$ws.Range("C13").Value = 1.4603098939655379

$ws.Range("B14").Value = 54.500806316273263
$ws.Range("C14").Value = 54.547011487965577

$ws.Range("B15").Value = 53.03289524862447
$ws.Range("C15").Value = 53.025335621785025

$ws.Range("C17").Value = -0.17012947120127564
$ws.Range("C18").Value = 0.0042909353465846326
$ws.Range("C19").Value = 0.0048958564920162737
$ws.Range("C20").Value = -1.5298114510016303
$ws.Range("C21").Value = -0.012826385019179898

$ws.Range("B22").Value = 8480
$ws.Range("C22").Value = 10352
